$d = $word.ActiveDocument

# --- Change 1: merge the "-- " lead-in run with the sentence that follows it,
# for the three CONCAT_WS() bullet paragraphs. Re-typing the text via Find &
# Replace (same formatting throughout) causes Word to collapse the range back
# down to a single run, matching the target XML shape.
$mergeTargets = @(
  "The first argument is the separator for the rest of the arguments.",
  "The separator is added between the strings to be concatenated.",
  "The separator can be a string, as can the rest of the arguments. If the separator is NULL, the result is NULL."
)
foreach ($t in $mergeTargets) {
  $find = $d.Content.Find
  $find.ClearFormatting()
  $find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

# --- Change 2: append the new "Lower(str)" / "Lpad(str, len, padstr)" section
# (plus the trailing LCASE/LEFT/CHAR_LENGTH/LENGTH heading stubs) right after
# the CONCAT_WS department example, before the existing trailing blank lines.
# A fresh blank paragraph is created first so InsertXML has a disposable
# paragraph to expand into, instead of overwriting the department example.
$anchor = $d.Paragraphs.Item(22)
$anchor.Range.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Item(23)
$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>Lower(str)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">-- </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Returns the string str with all characters changed to lowercase according to the current character set mapping.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>SELECT LOWER(''GDWLKGHDSLGHDG'');</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t>-- o/</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>p :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> "gdwlkghdslghdg"</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>LOWER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>first_name) as full_name FROM dev_schema.employee;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>LOWER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>CONCAT(first_name, '' '', last_name)) as full_name FROM dev_schema.employee;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>Lpad(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>str, len, padstr)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">-- </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>This string function adds the symbol or the string to the left side of the original string. In Structured Query Language, we can use this function on both string and numbers.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>LCASE(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>str)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>LEFT(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>str, len)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>CHAR_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>LENGTH(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>str)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>LENGTH(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>str)</w:t></w:r></w:p>'
$placeholder.Range.InsertXML($newContentXml)
